$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.235.21"
$ws.Range("E2").Value = "  -2.35%  "

$ws.Range("D3").Value = "1.852.11"
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("D4").Formula = "'" + "1.000"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Formula = "'" + "0.6976"
$ws.Range("E5").Value = "  -5.77%  "

$ws.Range("D6").Formula = "'" + "238.76"
$ws.Range("E6").Value = "  -1.59%  "

$ws.Range("D7").Formula = "'" + "1.000"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Formula = "'" + "0.3070"
$ws.Range("E8").Value = "  -2.35%  "

$ws.Range("D9").Formula = "'" + "0.07535"
$ws.Range("E9").Value = "  +4.52%  "

$ws.Range("D10").Formula = "'" + "23.61"
$ws.Range("E10").Value = "  -4.11%  "

$ws.Range("E11").Value = "  -2.91%  "

$ws.Range("D12").Formula = "'" + "0.7264"
$ws.Range("E12").Value = "  -3.38%  "

$ws.Range("D13").Value = "1.831.18"
$ws.Range("E13").Value = "  -3.09%  "

$ws.Range("D14").Formula = "'" + "5.190"
$ws.Range("E14").Value = "  -4.17%  "

$ws.Range("D15").Formula = "'" + "89.20"
$ws.Range("E15").Value = "  -3.70%  "

$ws.Range("D16").Value = "29.280.51"
$ws.Range("E16").Value = "  -2.25%  "

$ws.Range("D17").Formula = "'" + "5.867"
$ws.Range("E17").Value = "  -4.00%  "

$ws.Range("D18").Formula = "'" + "241.93"
$ws.Range("E18").Value = "  -3.12%  "

$ws.Range("D19").Formula = "'" + "0.000007714"
$ws.Range("E19").Value = "  -1.77%  "

$ws.Range("E20").Value = "  -3.45%  "

$ws.Range("D21").Formula = "'" + "1.001"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.101.60"
$ws.Range("E22").Value = "  -1.96%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Formula = "'" + "1.001"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Formula = "'" + "7.612"
$ws.Range("E24").Value = "  -5.27%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Formula = "'" + "9.051"
$ws.Range("E25").Value = "  -2.25%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Formula = "'" + "162.23"
$ws.Range("E26").Value = "  -1.65%  "

$ws.Range("D27").Formula = "'" + "0.1461"
$ws.Range("E27").Value = "  -5.83%  "

$ws.Range("D28").Formula = "'" + "18.07"
$ws.Range("E28").Value = "  -3.36%  "

$ws.Range("D29").Formula = "'" + "1.927"
$ws.Range("E29").Value = "  -5.26%  "

$ws.Range("D30").Formula = "'" + "1.399"
$ws.Range("E30").Value = "  -7.50%  "

$ws.Range("D31").Formula = "'" + "1.503"
$ws.Range("E31").Value = "  -2.07%  "

$ws.Range("D32").Formula = "'" + "4.425"
$ws.Range("E32").Value = "  -3.90%  "

$ws.Range("D33").Formula = "'" + "4.034"
$ws.Range("E33").Value = "  -5.92%  "

$ws.Range("E34").Value = "  -0.89%  "

$ws.Range("E35").Value = "  -3.49%  "

$ws.Range("D36").Formula = "'" + "0.7097"
$ws.Range("E36").Value = "  -5.02%  "

$ws.Range("D37").Formula = "'" + "0.9998"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("D38").Formula = "'" + "2.662"
$ws.Range("E38").Value = "  -1.34%  "

$ws.Range("D39").Formula = "'" + "0.01859"
$ws.Range("E39").Value = "  -5.58%  "

$ws.Range("D40").Formula = "'" + "2.704"
$ws.Range("E40").Value = "  -1.92%  "

$ws.Range("D41").Formula = "'" + "0.9273"
$ws.Range("E41").Value = "  +8.20%  "

$ws.Range("D42").Formula = "'" + "0.4300"
$ws.Range("E42").Value = "  -5.63%  "

$ws.Range("D43").Formula = "'" + "5.907"
$ws.Range("E43").Value = "  -3.82%  "

$ws.Range("D44").Value = "1.044.03"
$ws.Range("E44").Value = "  -6.07%  "

$ws.Range("D45").Formula = "'" + "69.51"
$ws.Range("E45").Value = "  -3.92%  "

$ws.Range("D46").Formula = "'" + "1.000"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").Formula = "'" + "102.45"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("D48").Formula = "'" + "7.220"
$ws.Range("E48").Value = "  -5.10%  "

$ws.Range("D49").Formula = "'" + "1.735"
$ws.Range("E49").Value = "  -6.42%  "

$ws.Range("D50").Formula = "'" + "9.233"
$ws.Range("E50").Value = "  -2.84%  "

$ws.Range("D51").Value = "1.987.03"
$ws.Range("E51").Value = "  -2.56%  "
